$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-format the "% Testing" cells for the Header sub-tasks (H5:H10).
#    These rows don't have a Testing phase, so their % cell should look
#    like the other "disabled / n.a." testing cells (same style already
#    used by H11 and the rest of the sheet: grey fill, percent format).
#    Copying the format from H11 reproduces Excel's own style re-use
#    (no new style entry gets created).
# ---------------------------------------------------------------------
$ws.Range("H11").Copy() | Out-Null
$ws.Range("H5:H10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Row 17 ("Responsive" under "About"): the Testing start/finish dates
#    were entered by mistake (copy/paste leftover) - clear them back out.
# ---------------------------------------------------------------------
$ws.Range("F17").ClearContents() | Out-Null
$ws.Range("G17").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# 3) "News" task (row 22) and its two sub-tasks "Grid Layout" (row 23)
#    and "Hover" (row 24) now have a Build start date, and the sub-tasks
#    are fully built (100%) which brings the parent task to 67%.
# ---------------------------------------------------------------------

# Donor cells already carrying the exact desired formats.
$ws.Range("B4").Copy() | Out-Null     # task-row "Start Build" date style
$ws.Range("B22").PasteSpecial(-4122) | Out-Null

$ws.Range("B5").Copy() | Out-Null     # sub-task-row "Start Build" date style
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null

$ws.Range("C5").Copy() | Out-Null     # sub-task-row "Finish Build" date style
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null

$ws.Range("B22").Value2 = 44473
$ws.Range("B23").Value2 = 44473
$ws.Range("C23").Value2 = 44474
$ws.Range("B24").Value2 = 44473
$ws.Range("C24").Value2 = 44474

$ws.Range("D22").Value2 = 0.67
$ws.Range("D23").Value2 = 1
$ws.Range("D24").Value2 = 1

# ---------------------------------------------------------------------
# 4) Scroll / selection housekeeping to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("J6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16

"done"
